$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (appears on the Overview sheet, columns B & C, rows 2 & 3)
# ---------------------------------------------------------------------------
foreach ($sheet in $wb.Worksheets) {
    $sheet.Cells.Replace("Ready for handoff", "Handed back: in sync with en-US")
}

# ---------------------------------------------------------------------------
# 2. Placeholder handback datetime "0001-01-01 00:00:00" becomes a real
#    timestamp. zh-cn's handback finished at 13:13:56, de-de's at 13:14:04.
# ---------------------------------------------------------------------------
foreach ($sheet in $wb.Worksheets) {
    $sheet.Cells.Replace("0001-01-01 00:00:00", "2016-03-23 13:13:56")
}

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-03-23 13:14:04"
$dede.Range("H3").Value = "2016-03-23 13:14:04"

# ---------------------------------------------------------------------------
# 3. Populate the "Latest Handback File" (F) / "Latest Handback DateTime"
#    hyperlink-bearing (G) columns for both the zh-cn and de-de reports, and
#    wire up the corresponding hyperlinks (same targets as the already
#    present handoff-file / handoff-target-file hyperlinks in columns A/D).
# ---------------------------------------------------------------------------
function Add-HandbackColumns($sheetName, $xlfDisplay) {
    $ws = $wb.Worksheets.Item($sheetName)

    $mdDisplay = "1b1eeac6-071f-4ba3-aff0-ee88e91b9af4.md"

    $mdUrl = $null
    $xlfUrl = $null
    foreach ($h in $ws.Hyperlinks) {
        $addr = $h.Range.Address()
        if ($addr -eq '$A$2') { $mdUrl = $h.Address }
        if ($addr -eq '$D$2') { $xlfUrl = $h.Address }
    }

    foreach ($row in 2, 3) {
        $ws.Range("F$row").Value = $mdDisplay
        $ws.Range("G$row").Value = $xlfDisplay

        $ws.Hyperlinks.Add($ws.Range("F$row"), $mdUrl, "", "", $mdDisplay) | Out-Null
        $ws.Hyperlinks.Add($ws.Range("G$row"), $xlfUrl, "", "", $xlfDisplay) | Out-Null
    }
}

Add-HandbackColumns "zh-cn" "1b1eeac6-071f-4ba3-aff0-ee88e91b9af4.895861227040b6a7f96e604f9220db9b47fde168.zh-cn.xlf"
Add-HandbackColumns "de-de" "1b1eeac6-071f-4ba3-aff0-ee88e91b9af4.895861227040b6a7f96e604f9220db9b47fde168.de-de.xlf"
